$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") for rows 2..535 from 45192 to 45202
$ws.Range("C2:C535").Value = 45202

# 2. Row 535 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(535).RowHeight = 15

# 3. Append new row 536 with the new logging notice
$ws.Range("A536").Value = "A 46437-2023"
$ws.Range("B536").Value = 45197
$ws.Range("C536").Value = 45202
$ws.Range("D536").Value = "ÖSTERGÖTLANDS LÄN"
$ws.Range("E536").Value = "MOTALA"
$ws.Range("F536").Value = "Holmen skog AB"
$ws.Range("G536").Value = 0.8
$ws.Range("H536").Value = 0
$ws.Range("I536").Value = 0
$ws.Range("J536").Value = 0
$ws.Range("K536").Value = 0
$ws.Range("L536").Value = 0
$ws.Range("M536").Value = 0
$ws.Range("N536").Value = 0
$ws.Range("O536").Value = 0
$ws.Range("P536").Value = 0
$ws.Range("Q536").Value = 0
$ws.Range("R536").Value = ""

$ws.Range("B536:C536").NumberFormat = "YYYY-MM-DD"
$ws.Range("R536").WrapText = $true
